$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A got wider (bug fix needed more room for the Date column) - closest
# reachable width to the target 14.85546875 "character" width.
$ws.Columns.Item(1).ColumnWidth = 14

# Row 4: trade that was missing a Date, a SellPrice, and was mis-flagged as
# not Profitable. Copy G1's date/number format (style index 1) onto A4 so the
# new Date value renders like the rest of the sheet, then write the value.
$ws.Range("G1").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = Get-Date -Year 2016 -Month 9 -Day 22 -Hour 15 -Minute 25 -Second 58

# Fill in the SellPrice that was missing.
$ws.Range("D4").Value = 82.03

# Mark the trade as Profitable (it was previously left blank / falsy).
$ws.Range("G4").Value = $true
